# Edit script: add "metadata" worksheet + refresh "time_taken"/query timestamps
# on the existing "data" worksheet (per commit "Refined metadata to be additional tab").

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Refresh the per-row query timestamps in column F of the "data" sheet.
#    (F2:F72 each get a new "panel_query_time"-style timestamp value.)
# ---------------------------------------------------------------------------
$timestamps = @(
  "2021-10-05 14:33:47.869786",
  "2021-10-05 14:33:47.869794",
  "2021-10-05 14:33:47.869797",
  "2021-10-05 14:33:47.869800",
  "2021-10-05 14:33:47.869803",
  "2021-10-05 14:33:47.869806",
  "2021-10-05 14:33:47.869808",
  "2021-10-05 14:33:47.869811",
  "2021-10-05 14:33:47.869814",
  "2021-10-05 14:33:47.869816",
  "2021-10-05 14:33:47.869819",
  "2021-10-05 14:33:47.869821",
  "2021-10-05 14:33:47.869824",
  "2021-10-05 14:33:47.869827",
  "2021-10-05 14:33:47.869829",
  "2021-10-05 14:33:47.869832",
  "2021-10-05 14:33:47.869834",
  "2021-10-05 14:33:47.869837",
  "2021-10-05 14:33:47.869839",
  "2021-10-05 14:33:47.869842",
  "2021-10-05 14:33:47.869844",
  "2021-10-05 14:33:47.869847",
  "2021-10-05 14:33:47.869849",
  "2021-10-05 14:33:47.869852",
  "2021-10-05 14:33:47.869854",
  "2021-10-05 14:33:47.869857",
  "2021-10-05 14:33:47.869859",
  "2021-10-05 14:33:47.869862",
  "2021-10-05 14:33:47.869864",
  "2021-10-05 14:33:47.869867",
  "2021-10-05 14:33:47.869876",
  "2021-10-05 14:33:47.869879",
  "2021-10-05 14:33:47.869882",
  "2021-10-05 14:33:47.869885",
  "2021-10-05 14:33:47.869887",
  "2021-10-05 14:33:47.869890",
  "2021-10-05 14:33:47.869892",
  "2021-10-05 14:33:47.869895",
  "2021-10-05 14:33:47.869897",
  "2021-10-05 14:33:47.869900",
  "2021-10-05 14:33:47.869903",
  "2021-10-05 14:33:47.869905",
  "2021-10-05 14:33:47.869908",
  "2021-10-05 14:33:47.869910",
  "2021-10-05 14:33:47.869913",
  "2021-10-05 14:33:47.869915",
  "2021-10-05 14:33:47.869918",
  "2021-10-05 14:33:47.869920",
  "2021-10-05 14:33:47.869923",
  "2021-10-05 14:33:47.869925",
  "2021-10-05 14:33:47.869928",
  "2021-10-05 14:33:47.869930",
  "2021-10-05 14:33:47.869933",
  "2021-10-05 14:33:47.869936",
  "2021-10-05 14:33:47.869938",
  "2021-10-05 14:33:47.869941",
  "2021-10-05 14:33:47.869943",
  "2021-10-05 14:33:47.869946",
  "2021-10-05 14:33:47.869948",
  "2021-10-05 14:33:47.869951",
  "2021-10-05 14:33:47.869953",
  "2021-10-05 14:33:47.869956",
  "2021-10-05 14:33:47.869958",
  "2021-10-05 14:33:47.869961",
  "2021-10-05 14:33:47.869965",
  "2021-10-05 14:33:47.869968",
  "2021-10-05 14:33:47.869970",
  "2021-10-05 14:33:47.869973",
  "2021-10-05 14:33:47.869975",
  "2021-10-05 14:33:47.869978",
  "2021-10-05 14:33:47.869980"
)
for ($i = 0; $i -lt $timestamps.Count; $i++) {
  $row = $i + 2
  $ws1.Range("F$row").Value = $timestamps[$i]
}

# ---------------------------------------------------------------------------
# 2. Add a new "metadata" worksheet right after "data".
# ---------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $ws1)
$newSheet.Name = "metadata"

# Copy the bold/bordered/centered header style used on the "data" sheet so the
# new header row (B1:G1) and the index cell (A2) match formatting.
$ws1.Range("B1:F1").Copy()
$newSheet.Range("B1:F1").PasteSpecial(-4122)
$ws1.Range("A2").Copy()
$newSheet.Range("G1").PasteSpecial(-4122)
$newSheet.Range("A2").PasteSpecial(-4122)

# Header row
$newSheet.Range("B1").Value = "data_name"
$newSheet.Range("C1").Value = "data_id"
$newSheet.Range("D1").Value = "data_version"
$newSheet.Range("E1").Value = "data_version_created"
$newSheet.Range("F1").Value = "panel_query_time"
$newSheet.Range("G1").Value = "panel_get_request"

# Data row
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "Early-onset Parkinson disease"
$newSheet.Range("C2").Value = 26
$newSheet.Range("D2").NumberFormat = "@"
$newSheet.Range("D2").Value = "0.120"
$newSheet.Range("D2").Style = "Normal"
$newSheet.Range("E2").Value = "2021-08-30T09:00:46.190463Z"
$newSheet.Range("F2").Value = "2021-10-05 14:33:47.866027"
$newSheet.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/26/?format=json"

# Keep "data" as the active sheet/tab, matching the original workbook view.
$ws1.Activate()
